$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1266-1267, shifting existing rows 1266:1326 down to 1268:1328
$ws.Range("A1266:A1267").EntireRow.Insert()

# Fill in the new row 1266 with the first new weekly data point
$ws.Range("A1266").Value = 10
$ws.Range("B1266").Value = "Vega Modelo de Temuco"
$ws.Range("C1266").Value = "La Araucanía"
$ws.Range("D1266").Value = 45267
$ws.Range("E1266").Value = 9
$ws.Range("F1266").Value = 100114001
$ws.Range("G1266").Value = "Papa"
$ws.Range("H1266").Value = "Asterix"
$ws.Range("I1266").Value = "1a (guarda)"
$ws.Range("J1266").Value = 180
$ws.Range("K1266").Value = 18000
$ws.Range("L1266").Value = 18000
$ws.Range("M1266").Value = 18000
$ws.Range("N1266").Value = "$/malla 25 kilos"
$ws.Range("O1266").Value = "Provincia de Cautín"
$ws.Range("P1266").Value = 720
$ws.Range("Q1266").Value = 25
$ws.Range("R1266").Value = "Hortaliza"

# Fill in the new row 1267 with the second new weekly data point
$ws.Range("A1267").Value = 10
$ws.Range("B1267").Value = "Vega Modelo de Temuco"
$ws.Range("C1267").Value = "La Araucanía"
$ws.Range("D1267").Value = 45267
$ws.Range("E1267").Value = 9
$ws.Range("F1267").Value = 100114001
$ws.Range("G1267").Value = "Papa"
$ws.Range("H1267").Value = "Cornado"
$ws.Range("I1267").Value = "1a nueva(o)"
$ws.Range("J1267").Value = 1200
$ws.Range("K1267").Value = 22000
$ws.Range("L1267").Value = 23000
$ws.Range("M1267").Value = 22583
$ws.Range("N1267").Value = "$/saco 25 kilos"
$ws.Range("O1267").Value = "Provincia de Cautín"
$ws.Range("P1267").Value = 903
$ws.Range("Q1267").Value = 25
$ws.Range("R1267").Value = "Hortaliza"
